$wb = $excel.ActiveWorkbook
$wsInventory = $wb.Worksheets.Item("Inventory")
$wsSales = $wb.Worksheets.Item("Sales")

# --- Inventory sheet: reduce on-hand quantities after sales ---
$wsInventory.Range("M2").Value = 15
$wsInventory.Range("M4").Value = 21
$wsInventory.Range("M5").Value = 201

# --- Sales sheet: append newly recorded sales transactions ---
# New rows, in order, with the same shape as the existing table
# (Item Desc, Price, Quantity, Total, Customer, Time)
$newSales = @(
    @("glplpp", 10.5,  2,  21,     "",       "2023-01-05 20:43:09"),
    @("lop",    10.5,  1,  10.5,   "selome", "2023-01-05 20:43:59"),
    @("lop",    10.5,  2,  21,     "selome", "2023-01-07 17:06:15"),
    @("lop",    10.5,  2,  21,     "selome", "2023-01-07 17:09:23"),
    @("lop",    10.5,  2,  21,     "selome", "2023-01-07 17:10:01"),
    @("lop",    10.5,  1,  10.5,   "selome", "2023-01-07 17:10:15"),
    @("emp",    12.35, 23, 284.05, "lewi",   "2023-01-07 17:12:15"),
    @("emp",    12.35, 23, 284.05, "lewi",   "2023-01-07 17:13:38"),
    @("emp",    12.35, 2,  24.7,   "Bogale", "2023-01-08 19:15:13")
)

$lastRow = 6
foreach ($sale in $newSales) {
    $newRowIndex = $lastRow + 1

    # Insert a new row, copying the formatting of the last existing row so
    # the appended row keeps the table's look (matches how extending the
    # Sales log in Excel behaves)
    $wsSales.Rows.Item($lastRow).Copy()
    $wsSales.Rows.Item($newRowIndex).Insert()

    $wsSales.Cells.Item($newRowIndex, 1).Value = $sale[0]
    $wsSales.Cells.Item($newRowIndex, 2).Value = $sale[1]
    $wsSales.Cells.Item($newRowIndex, 3).Value = $sale[2]
    $wsSales.Cells.Item($newRowIndex, 4).Value = $sale[3]
    $wsSales.Cells.Item($newRowIndex, 5).Value = $sale[4]
    $wsSales.Cells.Item($newRowIndex, 6).Value = $sale[5]

    $lastRow = $newRowIndex
}
